# Build site at 2022-09-26 16:07:08 UTC
# LOB1214.xlsx: remove the stray "Docentes responsaveis" data row (it held the
# professor's name with no label) and fix up several cells whose contents had
# been shifted/mixed up.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous row 13 only carried B13/C13 = "9146830 - Danubia Caporusso
# Bargos" with no A13 label and no row formatting. Delete it outright; Excel
# shifts every row below it up by one and keeps their formatting intact.
$ws.Rows(13).Delete()

# After the shift, a handful of B/C cells need their text content corrected.

# Row 10 ("Objetivos:") now shows the responsible professor instead of the
# long objectives paragraph.
$ws.Range("B10").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Range("C10").Value = "9146830 - Danúbia Caporusso Bargos"

# Row 13 ("Programa resumido:") now just says "Semestral".
$ws.Range("B13").Value = "Semestral"
$ws.Range("C13").Value = "Semestral"

# Row 15 ("Programa:") now holds the activation date. Build it through a
# TEXT() formula and paste back as a value so Excel stores it as a plain
# shared string instead of auto-converting the dd/mm/yyyy text into a date
# serial number (which would also allocate a brand-new cell style).
$ws.Range("B15").Formula = "=TEXT(""01/01/2018"",""@"")"
$ws.Range("B15").Copy()
$ws.Range("B15").PasteSpecial(-4163)

$ws.Range("C15").Formula = "=TEXT(""01/01/2018"",""@"")"
$ws.Range("C15").Copy()
$ws.Range("C15").PasteSpecial(-4163)

# Row 18 ("Método:") now shows the responsible professor instead of the
# teaching-method description.
$ws.Range("B18").Value = "9146830 - Danúbia Caporusso Bargos"
$ws.Range("C18").Value = "9146830 - Danúbia Caporusso Bargos"

# Row 19 ("Critério:") now carries the teaching-method description text.
$ws.Range("B19").Value = "Aulas teóricas e práticas, visitas técnicas e exercícios dirigidos. Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios."
$ws.Range("C19").Value = "Aulas teóricas e práticas, visitas técnicas e exercícios dirigidos. Avaliação baseada em provas, exercícios e trabalhos práticos e relatórios."

# Row 20 ("Norma de recuperação:") now carries the grading-criteria text.
$ws.Range("B20").Value = "Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios."
$ws.Range("C20").Value = "Média ponderada das notas atribuídas às provas, exercícios e trabalhos práticos e relatórios."

# Row 21 ("Bibliografia:") now carries the recovery-exam text instead of the
# full bibliography.
$ws.Range("B21").Value = "Provas e/ou exercícios dirigidos."
$ws.Range("C21").Value = "Provas e/ou exercícios dirigidos."
